$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("airline_financials")

# Row 186: 2025 Q2 - DAL
$ws.Cells.Item(186, 1).Value = 2025
$ws.Cells.Item(186, 2).Value = 2
$ws.Cells.Item(186, 3).Value = "DAL"
$ws.Cells.Item(186, 4).Value = 16648000000
$ws.Cells.Item(186, 5).Value = 13867000000
$ws.Cells.Item(186, 6).Value = 14546000000
$ws.Cells.Item(186, 7).Value = 2130000000
$ws.Cells.Item(186, 8).Value = 66417000000
$ws.Cells.Item(186, 9).Value = 77645000000
$ws.Cells.Item(186, 11).Value = 470000000

# Row 187: 2025 Q2 - UAL
$ws.Cells.Item(187, 1).Value = 2025
$ws.Cells.Item(187, 2).Value = 2
$ws.Cells.Item(187, 3).Value = "UAL"
$ws.Cells.Item(187, 4).Value = 15236000000
$ws.Cells.Item(187, 5).Value = 13836000000
$ws.Cells.Item(187, 6).Value = 13911000000
$ws.Cells.Item(187, 7).Value = 973000000
$ws.Cells.Item(187, 8).Value = 70088000000
$ws.Cells.Item(187, 9).Value = 84347000000
$ws.Cells.Item(187, 11).Value = 188000000

# Row 188: 2025 Q2 - AAL (values not yet reported)
$ws.Cells.Item(188, 1).Value = 2025
$ws.Cells.Item(188, 2).Value = 2
$ws.Cells.Item(188, 3).Value = "AAL"

# Row 189: 2025 Q2 - LUV (values not yet reported)
$ws.Cells.Item(189, 1).Value = 2025
$ws.Cells.Item(189, 2).Value = 2
$ws.Cells.Item(189, 3).Value = "LUV"

# Restore selection state to match the saved workbook
$ws.Range("F192").Select() | Out-Null
